# Update "Ventas objetivo" related recalculated columns (L, R, T, U) and the
# summary metrics (C56, C67) on sheet "Semana_7" to reflect the new formula
# used for "Ventas objetivo".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

# Row 5
$ws.Range("L5").Value = 0
$ws.Range("R5").Value = 106
$ws.Range("T5").Value = 0

# Row 8
$ws.Range("L8").Value = 0

# Row 10
$ws.Range("L10").Value = 0
$ws.Range("R10").Value = 1
$ws.Range("T10").Value = 3

# Row 13
$ws.Range("R13").Value = 1

# Row 16
$ws.Range("L16").Value = 0

# Row 19
$ws.Range("L19").Value = 0

# Row 23
$ws.Range("L23").Value = 0
$ws.Range("R23").Value = 2
$ws.Range("T23").Value = 1
$ws.Range("U23").Value = 8

# Row 24
$ws.Range("R24").Value = 5
$ws.Range("T24").Value = 0

# Row 25
$ws.Range("L25").Value = 0
$ws.Range("R25").Value = 9
$ws.Range("T25").Value = 0
$ws.Range("U25").Value = 3

# Row 26
$ws.Range("L26").Value = 0

# Row 30
$ws.Range("L30").Value = 0

# Row 31
$ws.Range("L31").Value = 0
$ws.Range("R31").Value = 1
$ws.Range("T31").Value = 0

# Row 33
$ws.Range("L33").Value = 0
$ws.Range("R33").Value = 4
$ws.Range("T33").Value = 7
$ws.Range("U33").Value = 10

# Row 37
$ws.Range("R37").Value = 2

# Row 39
$ws.Range("L39").Value = 0
$ws.Range("R39").Value = 4

# Row 40
$ws.Range("L40").Value = 0

# Row 41
$ws.Range("L41").Value = 0
$ws.Range("R41").Value = 4
$ws.Range("T41").Value = 10
$ws.Range("U41").Value = 7

# Row 44
$ws.Range("R44").Value = 1

# Row 45
$ws.Range("L45").Value = 0
$ws.Range("R45").Value = 3
$ws.Range("T45").Value = 1
$ws.Range("U45").Value = 7

# Row 46
$ws.Range("L46").Value = 0
$ws.Range("R46").Value = 2

# Row 48
$ws.Range("L48").Value = 0
$ws.Range("R48").Value = 3

# Row 49
$ws.Range("L49").Value = 0
$ws.Range("R49").Value = 11
$ws.Range("T49").Value = 0

# Row 50
$ws.Range("R50").Value = 2
$ws.Range("T50").Value = 0

# Row 51
$ws.Range("L51").Value = 0
$ws.Range("R51").Value = 3

# Row 52
$ws.Range("L52").Value = 0

# Row 53
$ws.Range("R53").Value = 1
$ws.Range("T53").Value = 4
$ws.Range("U53").Value = 2

# Summary metrics
$ws.Range("C56").Value = 124
$ws.Range("C67").Value = 0
